# Refresh the cryptos price list (Price / Volume(1h) columns), matching the
# GitHub Actions scrape commit. Two coin pairs (OKB/BinanceUSD and
# ImmutableX/InternetComputer) also swapped table position.
#
# Price values that look like plain numbers (e.g. "1.002") are written with
# a leading apostrophe so Excel stores them as text (matching the source
# data, which keeps prices as literal strings such as "24.780.01" that are
# not valid numbers) instead of silently coercing them to numeric cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.780.01"
$ws.Range("E2").Value = "  +0.15%  "

$ws.Range("D3").Value = "1.705.37"
$ws.Range("E3").Value = "  +0.36%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'317.04"
$ws.Range("E5").Value = "  +0.12%  "

$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  -0.11%  "

$ws.Range("D7").Value = "'0.3931"
$ws.Range("E7").Value = "  -0.44%  "

$ws.Range("D8").Value = "'0.4048"
$ws.Range("E8").Value = "  -0.66%  "

$ws.Range("D9").Value = "'1.488"
$ws.Range("E9").Value = "  -1.26%  "

$ws.Range("B10").Value = "BinanceUSD"
$ws.Range("C10").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D10").Value = "'1.003"
$ws.Range("E10").Value = "  -0.14%  "

$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "'53.53"
$ws.Range("E11").Value = "  +1.32%  "

$ws.Range("D12").Value = "'0.08798"
$ws.Range("E12").Value = "  -1.17%  "

$ws.Range("D13").Value = "'26.35"
$ws.Range("E13").Value = "  +9.91%  "

$ws.Range("D14").Value = "'7.474"
$ws.Range("E14").Value = "  -3.77%  "

$ws.Range("D15").Value = "'8.093"
$ws.Range("E15").Value = "  -1.15%  "

$ws.Range("D16").Value = "'0.00001355"
$ws.Range("E16").Value = "  +2.05%  "

$ws.Range("D17").Value = "1.739.94"
$ws.Range("E17").Value = "  +1.81%  "

$ws.Range("D18").Value = "'96.24"
$ws.Range("E18").Value = "  -3.47%  "

$ws.Range("D19").Value = "'0.07188"
$ws.Range("E19").Value = "  +0.87%  "

$ws.Range("D20").Value = "'20.89"
$ws.Range("E20").Value = "  +4.67%  "

$ws.Range("D21").Value = "'7.304"
$ws.Range("E21").Value = "  +1.68%  "

$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  -0.53%  "

$ws.Range("D23").Value = "'14.35"
$ws.Range("E23").Value = "  -2.30%  "

$ws.Range("D24").Value = "24.782.00"
$ws.Range("E24").Value = "  +0.13%  "

$ws.Range("D25").Value = "'2.990"
$ws.Range("E25").Value = "  -4.35%  "

$ws.Range("D26").Value = "'2.347"
$ws.Range("E26").Value = "  +0.09%  "

$ws.Range("D27").Value = "'23.23"
$ws.Range("E27").Value = "  +0.41%  "

$ws.Range("D28").Value = "'166.61"
$ws.Range("E28").Value = "  +0.74%  "

$ws.Range("D29").Value = "'5.986"
$ws.Range("E29").Value = "  +15.59%  "

$ws.Range("D30").Value = "'8.492"
$ws.Range("E30").Value = "  -7.97%  "

$ws.Range("D31").Value = "'144.84"
$ws.Range("E31").Value = "  +4.27%  "

$ws.Range("D32").Value = "1.924.99"
$ws.Range("E32").Value = "  +1.49%  "

$ws.Range("D33").Value = "'2.268"
$ws.Range("E33").Value = "  +15.29%  "

$ws.Range("D34").Value = "'0.08790"
$ws.Range("E34").Value = "  -3.30%  "

$ws.Range("D35").Value = "'0.03139"
$ws.Range("E35").Value = "  +3.86%  "

$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "'7.191"
$ws.Range("E36").Value = "  -11.48%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.039"
$ws.Range("E37").Value = "  -3.72%  "

$ws.Range("D38").Value = "'0.2843"
$ws.Range("E38").Value = "  +0.53%  "

$ws.Range("D39").Value = "'10.86"
$ws.Range("E39").Value = "  -2.41%  "

$ws.Range("D40").Value = "'0.8298"
$ws.Range("E40").Value = "  +6.44%  "

$ws.Range("D41").Value = "'0.09216"
$ws.Range("E41").Value = "  -0.79%  "

$ws.Range("D42").Value = "'14.10"
$ws.Range("E42").Value = "  -2.25%  "

$ws.Range("D43").Value = "'1.473"
$ws.Range("E43").Value = "  -0.03%  "

$ws.Range("D44").Value = "'17.37"
$ws.Range("E44").Value = "  +7.51%  "

$ws.Range("D45").Value = "'2.682"
$ws.Range("E45").Value = "  +1.17%  "

$ws.Range("D46").Value = "'0.7375"
$ws.Range("E46").Value = "  +1.82%  "

$ws.Range("D47").Value = "'4.249"
$ws.Range("E47").Value = "  +0.31%  "

$ws.Range("D48").Value = "'1.390"
$ws.Range("E48").Value = "  +1.70%  "

$ws.Range("D49").Value = "'1.002"
$ws.Range("E49").Value = "  -0.11%  "

$ws.Range("D50").Value = "'140.70"
$ws.Range("E50").Value = "  +0.02%  "

$ws.Range("D51").Value = "'0.08272"
$ws.Range("E51").Value = "  +3.19%  "
